$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue "D2" "26.105.46"
Set-TextValue "E2" "  -0.95%  "
Set-TextValue "D3" "1.665.43"
Set-TextValue "E3" "  -1.32%  "
Set-TextValue "E4" "  -0.82%  "
Set-TextValue "D5" "209.34"
Set-TextValue "E5" "  -4.22%  "
Set-TextValue "D6" "0.5170"
Set-TextValue "E6" "  -4.75%  "
Set-TextValue "E7" "  -0.84%  "
Set-TextValue "D8" "0.2629"
Set-TextValue "E8" "  -3.85%  "
Set-TextValue "D9" "0.06202"
Set-TextValue "E9" "  -3.76%  "
Set-TextValue "D10" "21.05"
Set-TextValue "E10" "  -4.03%  "
Set-TextValue "D11" "0.07488"
Set-TextValue "E11" "  -2.50%  "
Set-TextValue "D12" "1.675.97"
Set-TextValue "E12" "  -0.70%  "
Set-TextValue "D13" "4.400"
Set-TextValue "E13" "  -2.75%  "
Set-TextValue "D14" "0.5556"
Set-TextValue "E14" "  -4.37%  "
$ws.Range("B15").Value = "Litecoin"
$ws.Range("C15").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue "D15" "65.68"
Set-TextValue "E15" "  +1.05%  "
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue "D16" "0.000007839"
Set-TextValue "E16" "  -6.37%  "
Set-TextValue "D17" "26.118.53"
Set-TextValue "E17" "  -1.14%  "
Set-TextValue "E18" "  -0.86%  "
Set-TextValue "D19" "4.764"
Set-TextValue "E19" "  -3.28%  "
Set-TextValue "D20" "10.35"
Set-TextValue "E20" "  -5.53%  "
Set-TextValue "D21" "185.55"
Set-TextValue "E21" "  -2.93%  "
Set-TextValue "D22" "6.136"
Set-TextValue "E22" "  -1.76%  "
Set-TextValue "E23" "  -0.85%  "
Set-TextValue "D24" "147.43"
Set-TextValue "E24" "  -1.42%  "
Set-TextValue "D25" "0.1232"
Set-TextValue "E25" "  -6.88%  "
Set-TextValue "D26" "7.520"
Set-TextValue "E26" "  -4.23%  "
Set-TextValue "D27" "15.82"
Set-TextValue "E27" "  +0.79%  "
Set-TextValue "D28" "0.06253"
Set-TextValue "E28" "  -1.43%  "
Set-TextValue "D29" "1.360"
Set-TextValue "E29" "  -3.29%  "
Set-TextValue "D30" "1.271"
Set-TextValue "E30" "  -4.10%  "
Set-TextValue "E31" "  -2.72%  "
Set-TextValue "D32" "3.402"
Set-TextValue "E32" "  -5.04%  "
Set-TextValue "D33" "1.617"
Set-TextValue "E33" "  -3.87%  "
Set-TextValue "D34" "0.9925"
Set-TextValue "E34" "  -4.48%  "
Set-TextValue "E35" "  -0.30%  "
Set-TextValue "D36" "0.6005"
Set-TextValue "E36" "  -2.21%  "
Set-TextValue "D37" "2.697"
Set-TextValue "E37" "  -0.39%  "
Set-TextValue "D38" "6.110"
Set-TextValue "E38" "  -2.46%  "
Set-TextValue "D39" "0.01601"
Set-TextValue "E39" "  -1.55%  "
Set-TextValue "D40" "1.069.17"
Set-TextValue "E40" "  -3.73%  "
Set-TextValue "D41" "0.8611"
Set-TextValue "E41" "  -2.32%  "
Set-TextValue "D42" "1.003"
Set-TextValue "E42" "  -1.22%  "
Set-TextValue "D43" "98.99"
Set-TextValue "E43" "  -2.73%  "
Set-TextValue "D44" "1.812.41"
Set-TextValue "E44" "  -1.42%  "
Set-TextValue "D45" "0.00000000108"
Set-TextValue "E45" "  -1.73%  "
Set-TextValue "D46" "55.80"
Set-TextValue "E46" "  -2.71%  "
Set-TextValue "D47" "1.002"
Set-TextValue "E47" "  -0.87%  "
Set-TextValue "E48" "  -0.34%  "
Set-TextValue "D49" "7.929"
Set-TextValue "E49" "  -3.18%  "
Set-TextValue "E50" "  -1.26%  "
Set-TextValue "D51" "5.894"
Set-TextValue "E51" "  -2.32%  "
